$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update date (use leading apostrophe so Excel stores it as literal
# text instead of auto-converting the "DD-MMM-YY"-looking string to a date
# serial) and the fare figures.
$ws.Range("A2").Formula = "'15-JAN-26"
$ws.Range("D2").Value = 421
$ws.Range("E2").Value = 458
$ws.Range("F2").Value = -37

# Row 3 (05-MAR-26 / SM-322 / Nile Air NP-120 threat entry) was removed
# entirely, shrinking the used range from A1:K3 to A1:K2.
$ws.Rows("3").Delete()
